$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns E and F (old "y" and "alpha" columns).
# This shifts old G (co_x) -> E and old H (co_y) -> F, and shrinks
# the used range from A1:H3 down to A1:F3, matching the target diff.
$ws.Range("E1:F1").EntireColumn.Delete() | Out-Null

# Ensure the new header text values are exactly as intended.
$ws.Range("E1").Value = "co_x"
$ws.Range("F1").Value = "co_y"
